$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1739
$ws.Range("F4").Value = 141
$ws.Range("F5").Value = 392
$ws.Range("F6").Value = 777
$ws.Range("F8").Value = 1110
$ws.Range("F9").Value = 292
$ws.Range("F11").Value = 435
$ws.Range("F12").Value = 645
$ws.Range("F14").Value = 500
$ws.Range("F17").Value = 159
$ws.Range("F18").Value = 2850
$ws.Range("F19").Value = 2596
$ws.Range("F24").Value = 218
$ws.Range("F26").Value = 158
$ws.Range("F27").Value = 587
$ws.Range("F28").Value = 972
$ws.Range("F29").Value = 12
$ws.Range("F31").Value = 264
$ws.Range("F32").Value = 1050
$ws.Range("F33").Value = 71
$ws.Range("F34").Value = 45
$ws.Range("F35").Value = 276

# Sheet: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1060
$ws.Range("F5").Value = 1060
$ws.Range("F11").Value = 5
$ws.Range("F15").Value = 593
$ws.Range("F16").Value = 96
$ws.Range("F18").Value = 972
$ws.Range("F21").Value = 612
$ws.Range("F25").Value = 300
$ws.Range("F26").Value = 267
$ws.Range("F27").Value = 3814
$ws.Range("F32").Value = 32

# Sheet: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1771
$ws.Range("F5").Value = 2417
$ws.Range("F6").Value = 1003
$ws.Range("F9").Value = 1267
$ws.Range("F10").Value = 335

# Sheet: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1771
$ws.Range("F4").Value = 2417
$ws.Range("F5").Value = 1739
$ws.Range("F6").Value = 1003
$ws.Range("F7").Value = 1267
$ws.Range("F8").Value = 335
$ws.Range("F10").Value = 141
$ws.Range("F11").Value = 392
$ws.Range("F12").Value = 777
$ws.Range("F15").Value = 1110
$ws.Range("F16").Value = 292
$ws.Range("F17").Value = 437
$ws.Range("F18").Value = 645
$ws.Range("F19").Value = 1060
$ws.Range("F20").Value = 500
$ws.Range("F23").Value = 159
$ws.Range("F24").Value = 2850
$ws.Range("F25").Value = 2596
$ws.Range("F30").Value = 218
$ws.Range("F31").Value = 158
$ws.Range("F32").Value = 587
$ws.Range("F33").Value = 972
$ws.Range("F34").Value = 593
$ws.Range("F35").Value = 593
$ws.Range("F36").Value = 12
$ws.Range("F37").Value = 96
$ws.Range("F39").Value = 264
$ws.Range("F44").Value = 300
$ws.Range("F45").Value = 300
$ws.Range("F46").Value = 267
$ws.Range("F47").Value = 1050
$ws.Range("F50").Value = 45
$ws.Range("F51").Value = 276
